# Update the "general" sheet: add two new rows for prior-distribution settings
# (flux prior and thermodynamic prior), pushing the existing rows below down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert two new blank rows above the current row 6
# ("Number of exp. conditions...") to hold the new settings.
$ws.Rows.Item(6).Resize(2).Insert()

# Fill in the new rows with label/value pairs.
$ws.Cells.Item(6, 1).Value2 = "Prior distribution for fluxes (uniform or normal)"
$ws.Cells.Item(6, 2).Value2 = "normal"
$ws.Cells.Item(7, 1).Value2 = "Prior distribution for thermodynamic quantities (uniform or normal)"
$ws.Cells.Item(7, 2).Value2 = "normal"

# Style column A (labels) like the other label cells in the top block (A2:A5):
# bold, left/top aligned, thin border all around.
$labels = $ws.Range("A6:A7")
$labels.HorizontalAlignment = -4131   # xlLeft
$labels.VerticalAlignment = -4160     # xlTop
$labels.Font.Bold = $true
$labels.Borders.LineStyle = 1

# Style column B (values) distinctly: centered, bottom aligned, thin border,
# regular (non-bold) font - matching the new value style used elsewhere.
$values = $ws.Range("B6:B7")
$values.HorizontalAlignment = -4108   # xlCenter
$values.VerticalAlignment = -4107     # xlBottom
$values.Font.Bold = $false
$values.Borders.LineStyle = 1

# Match the row height used in the block above (rows 2-5).
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8

# Make "general" the active sheet, with the two new rows selected -
# matching the reviewed state after the edit.
$ws.Activate()
$ws.Range("A6:B7").Select()

$wb.Save()
